$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 59, pushing the existing row 59
# (Agricola del Norte / Guayaba, 2021-10-04) down to row 61 unchanged.
$ws.Range("A59:A60").EntireRow.Insert()

# Update rows 38-58 in place: these keep the same "shape" of record
# (same market/product/category/variety) but several rows now carry
# different date / volume / price values per the weekly refresh.
$ws.Cells.Item(38, 4).Value = 45068
$ws.Cells.Item(38, 12).Value = 'Primera'
$ws.Cells.Item(38, 13).Value = 70
$ws.Cells.Item(38, 14).Value = 8000
$ws.Cells.Item(38, 15).Value = 9000
$ws.Cells.Item(38, 16).Value = 8571
$ws.Cells.Item(38, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(38, 19).Value = 857
$ws.Cells.Item(38, 20).Value = 10

$ws.Cells.Item(39, 4).Value = 44764
$ws.Cells.Item(39, 12).Value = 'Primera'
$ws.Cells.Item(39, 13).Value = 200
$ws.Cells.Item(39, 14).Value = 500
$ws.Cells.Item(39, 15).Value = 600
$ws.Cells.Item(39, 16).Value = 550
$ws.Cells.Item(39, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(39, 19).Value = 550
$ws.Cells.Item(39, 20).Value = 1

$ws.Cells.Item(40, 4).Value = 44316
$ws.Cells.Item(40, 12).Value = 'Primera'
$ws.Cells.Item(40, 13).Value = 140
$ws.Cells.Item(40, 14).Value = 1100
$ws.Cells.Item(40, 15).Value = 1200
$ws.Cells.Item(40, 16).Value = 1150
$ws.Cells.Item(40, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(40, 19).Value = 1150
$ws.Cells.Item(40, 20).Value = 1

$ws.Cells.Item(41, 4).Value = 44750
$ws.Cells.Item(41, 12).Value = 'Primera'
$ws.Cells.Item(41, 13).Value = 200
$ws.Cells.Item(41, 14).Value = 700
$ws.Cells.Item(41, 15).Value = 800
$ws.Cells.Item(41, 16).Value = 750
$ws.Cells.Item(41, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(41, 19).Value = 750
$ws.Cells.Item(41, 20).Value = 1

$ws.Cells.Item(42, 4).Value = 45030
$ws.Cells.Item(42, 12).Value = 'Primera'
$ws.Cells.Item(42, 13).Value = 100
$ws.Cells.Item(42, 14).Value = 1900
$ws.Cells.Item(42, 15).Value = 2000
$ws.Cells.Item(42, 16).Value = 1950
$ws.Cells.Item(42, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(42, 19).Value = 1950
$ws.Cells.Item(42, 20).Value = 1

$ws.Cells.Item(43, 4).Value = 44767
$ws.Cells.Item(43, 12).Value = 'Segunda'
$ws.Cells.Item(43, 13).Value = 200
$ws.Cells.Item(43, 14).Value = 800
$ws.Cells.Item(43, 15).Value = 900
$ws.Cells.Item(43, 16).Value = 850
$ws.Cells.Item(43, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(43, 19).Value = 850
$ws.Cells.Item(43, 20).Value = 1

$ws.Cells.Item(44, 4).Value = 45016
$ws.Cells.Item(44, 12).Value = 'Primera'
$ws.Cells.Item(44, 13).Value = 90
$ws.Cells.Item(44, 14).Value = 2400
$ws.Cells.Item(44, 15).Value = 2500
$ws.Cells.Item(44, 16).Value = 2450
$ws.Cells.Item(44, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(44, 19).Value = 2450
$ws.Cells.Item(44, 20).Value = 1

$ws.Cells.Item(45, 4).Value = 44767
$ws.Cells.Item(45, 12).Value = 'Primera'
$ws.Cells.Item(45, 13).Value = 200
$ws.Cells.Item(45, 14).Value = 800
$ws.Cells.Item(45, 15).Value = 900
$ws.Cells.Item(45, 16).Value = 850
$ws.Cells.Item(45, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(45, 19).Value = 850
$ws.Cells.Item(45, 20).Value = 1

$ws.Cells.Item(46, 4).Value = 44767
$ws.Cells.Item(46, 12).Value = 'Segunda'
$ws.Cells.Item(46, 13).Value = 140
$ws.Cells.Item(46, 14).Value = 750
$ws.Cells.Item(46, 15).Value = 800
$ws.Cells.Item(46, 16).Value = 775
$ws.Cells.Item(46, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(46, 19).Value = 775
$ws.Cells.Item(46, 20).Value = 1

$ws.Cells.Item(47, 4).Value = 44350
$ws.Cells.Item(47, 12).Value = 'Primera'
$ws.Cells.Item(47, 13).Value = 140
$ws.Cells.Item(47, 14).Value = 750
$ws.Cells.Item(47, 15).Value = 800
$ws.Cells.Item(47, 16).Value = 775
$ws.Cells.Item(47, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(47, 19).Value = 775
$ws.Cells.Item(47, 20).Value = 1

$ws.Cells.Item(48, 4).Value = 44389
$ws.Cells.Item(48, 12).Value = 'Segunda'
$ws.Cells.Item(48, 13).Value = 120
$ws.Cells.Item(48, 14).Value = 600
$ws.Cells.Item(48, 15).Value = 700
$ws.Cells.Item(48, 16).Value = 650
$ws.Cells.Item(48, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(48, 19).Value = 650
$ws.Cells.Item(48, 20).Value = 1

$ws.Cells.Item(49, 4).Value = 44725
$ws.Cells.Item(49, 12).Value = 'Primera'
$ws.Cells.Item(49, 13).Value = 140
$ws.Cells.Item(49, 14).Value = 700
$ws.Cells.Item(49, 15).Value = 800
$ws.Cells.Item(49, 16).Value = 750
$ws.Cells.Item(49, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(49, 19).Value = 750
$ws.Cells.Item(49, 20).Value = 1

$ws.Cells.Item(50, 4).Value = 44725
$ws.Cells.Item(50, 12).Value = 'Segunda'
$ws.Cells.Item(50, 13).Value = 160
$ws.Cells.Item(50, 14).Value = 500
$ws.Cells.Item(50, 15).Value = 600
$ws.Cells.Item(50, 16).Value = 550
$ws.Cells.Item(50, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(50, 19).Value = 550
$ws.Cells.Item(50, 20).Value = 1

$ws.Cells.Item(51, 4).Value = 44690
$ws.Cells.Item(51, 12).Value = 'Primera'
$ws.Cells.Item(51, 13).Value = 100
$ws.Cells.Item(51, 14).Value = 1600
$ws.Cells.Item(51, 15).Value = 1700
$ws.Cells.Item(51, 16).Value = 1650
$ws.Cells.Item(51, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(51, 19).Value = 1650
$ws.Cells.Item(51, 20).Value = 1

$ws.Cells.Item(52, 4).Value = 44778
$ws.Cells.Item(52, 12).Value = 'Primera'
$ws.Cells.Item(52, 13).Value = 200
$ws.Cells.Item(52, 14).Value = 700
$ws.Cells.Item(52, 15).Value = 800
$ws.Cells.Item(52, 16).Value = 750
$ws.Cells.Item(52, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(52, 19).Value = 750
$ws.Cells.Item(52, 20).Value = 1

$ws.Cells.Item(53, 4).Value = 44778
$ws.Cells.Item(53, 12).Value = 'Segunda'
$ws.Cells.Item(53, 13).Value = 140
$ws.Cells.Item(53, 14).Value = 500
$ws.Cells.Item(53, 15).Value = 600
$ws.Cells.Item(53, 16).Value = 550
$ws.Cells.Item(53, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(53, 19).Value = 550
$ws.Cells.Item(53, 20).Value = 1

$ws.Cells.Item(54, 4).Value = 44417
$ws.Cells.Item(54, 12).Value = 'Primera'
$ws.Cells.Item(54, 13).Value = 200
$ws.Cells.Item(54, 14).Value = 1300
$ws.Cells.Item(54, 15).Value = 1400
$ws.Cells.Item(54, 16).Value = 1350
$ws.Cells.Item(54, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(54, 19).Value = 1350
$ws.Cells.Item(54, 20).Value = 1

$ws.Cells.Item(55, 4).Value = 44694
$ws.Cells.Item(55, 12).Value = 'Primera'
$ws.Cells.Item(55, 13).Value = 120
$ws.Cells.Item(55, 14).Value = 1400
$ws.Cells.Item(55, 15).Value = 1500
$ws.Cells.Item(55, 16).Value = 1450
$ws.Cells.Item(55, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(55, 19).Value = 1450
$ws.Cells.Item(55, 20).Value = 1

$ws.Cells.Item(56, 4).Value = 44694
$ws.Cells.Item(56, 12).Value = 'Segunda'
$ws.Cells.Item(56, 13).Value = 140
$ws.Cells.Item(56, 14).Value = 1100
$ws.Cells.Item(56, 15).Value = 1200
$ws.Cells.Item(56, 16).Value = 1150
$ws.Cells.Item(56, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(56, 19).Value = 1150
$ws.Cells.Item(56, 20).Value = 1

$ws.Cells.Item(57, 4).Value = 44729
$ws.Cells.Item(57, 12).Value = 'Primera'
$ws.Cells.Item(57, 13).Value = 150
$ws.Cells.Item(57, 14).Value = 700
$ws.Cells.Item(57, 15).Value = 800
$ws.Cells.Item(57, 16).Value = 750
$ws.Cells.Item(57, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(57, 19).Value = 750
$ws.Cells.Item(57, 20).Value = 1

$ws.Cells.Item(58, 4).Value = 44729
$ws.Cells.Item(58, 12).Value = 'Segunda'
$ws.Cells.Item(58, 13).Value = 160
$ws.Cells.Item(58, 14).Value = 500
$ws.Cells.Item(58, 15).Value = 600
$ws.Cells.Item(58, 16).Value = 550
$ws.Cells.Item(58, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(58, 19).Value = 550
$ws.Cells.Item(58, 20).Value = 1

# Fill the two newly inserted rows (59-60) with new records for the
# same market/product/category/variety (Primera & Segunda quality).
$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(59, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(59, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(59, 4).Value = 44722
$ws.Cells.Item(59, 5).Value = 15
$ws.Cells.Item(59, 6).Value = 'Fruta'
$ws.Cells.Item(59, 7).Value = 100108
$ws.Cells.Item(59, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(59, 9).Value = 100108001
$ws.Cells.Item(59, 10).Value = 'Guayaba'
$ws.Cells.Item(59, 11).Value = 'Sin especificar'
$ws.Cells.Item(59, 12).Value = 'Primera'
$ws.Cells.Item(59, 13).Value = 140
$ws.Cells.Item(59, 14).Value = 800
$ws.Cells.Item(59, 15).Value = 900
$ws.Cells.Item(59, 16).Value = 850
$ws.Cells.Item(59, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(59, 18).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(59, 19).Value = 850
$ws.Cells.Item(59, 20).Value = 1

$ws.Cells.Item(60, 1).Value = 1
$ws.Cells.Item(60, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(60, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(60, 4).Value = 44722
$ws.Cells.Item(60, 5).Value = 15
$ws.Cells.Item(60, 6).Value = 'Fruta'
$ws.Cells.Item(60, 7).Value = 100108
$ws.Cells.Item(60, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(60, 9).Value = 100108001
$ws.Cells.Item(60, 10).Value = 'Guayaba'
$ws.Cells.Item(60, 11).Value = 'Sin especificar'
$ws.Cells.Item(60, 12).Value = 'Segunda'
$ws.Cells.Item(60, 13).Value = 200
$ws.Cells.Item(60, 14).Value = 700
$ws.Cells.Item(60, 15).Value = 800
$ws.Cells.Item(60, 16).Value = 750
$ws.Cells.Item(60, 17).Value = '$/kilo (en caja de 10 kilos )'
$ws.Cells.Item(60, 18).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(60, 19).Value = 750
$ws.Cells.Item(60, 20).Value = 1

Write-Host "Done. Dimension:" $ws.UsedRange.Address